# Apply updated crypto price/volume data as per commit "Updated cryptos list".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.308.56'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.078.29'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.00'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.667'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.12'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +24.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '63.31'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.386'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0761'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.106'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +6.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.41'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.381.63'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.843'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.31'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.075.54'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.206.09'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.24'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.65'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +12.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0867'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.55%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.29'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.61%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.48'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '170.95'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.59%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.17'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.71'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.04'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.31%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.10'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +22.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '22.62'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.59'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0634'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.96%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0917'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -6.02%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.37'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +6.17%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.29'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'WEMIXToken'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.84'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.36'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.104'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +24.08%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '18.21'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +11.75%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0229'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.17'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '99.54'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.32'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +84.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.79'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.43'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +6.40%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.322.51'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.93'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.04'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.61%  '
